# Convert the " m:'doc.html'.fromHTMLURI() " field (begin/instrText/end)
# into plain literal text runs "{", "m", ":", "'", "doc.html", "'.fromHTMLURI()", "}"
# while keeping the _GoBack bookmark in place (TokenIteratorFieldRewriterSplit
# emits the M2Doc token delimiters as separate w:t runs instead of a Word field).

$d = $word.ActiveDocument

# Locate the paragraph that holds the M2Doc field (its code contains
# "fromHTMLURI") without assuming a fixed paragraph index.
$paraIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Fields.Count -gt 0) {
        foreach ($fld in $p.Range.Fields) {
            if ($fld.Code.Text -match "fromHTMLURI") {
                $paraIndex = $i
            }
        }
    }
}

if ($paraIndex -eq 0) {
    # Fallback: no paragraph matched by field code text, just use the
    # paragraph that holds the document's first field.
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($paraIndex -eq 0 -and $p.Range.Fields.Count -gt 0) {
            $paraIndex = $i
        }
    }
}

$targetPara = $d.Paragraphs.Item($paraIndex)
$targetField = $targetPara.Range.Fields.Item(1)

# Remove the field (begin/instrText.../end) but keep the (now empty)
# paragraph it lived in.
$targetField.Delete()

$p2 = $d.Paragraphs.Item($paraIndex)
$ins = $p2.Range
$ins.Collapse(1)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t>&apos;</w:t></w:r><w:r><w:t>doc.html</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>&apos;.fromHTMLURI()</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$ins.InsertXML($xmlFrag)

Write-Host "Rewrote field into literal-text runs in paragraph $paraIndex"
